$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "dam" subdomain to "photos" across the affected rows (14-17).
# Editing the cell values also updates the dependent formula results in
# column A (e.g. "dam.bcparks.ca" -> "photos.bcparks.ca") and the
# shared-strings table (the now-unused "dam" string is dropped and a new
# "photos" string is appended).
$ws.Range("B14").Value = "photos"
$ws.Range("B15").Value = "photos"
$ws.Range("B16").Value = "photos"
$ws.Range("B17").Value = "photos"

# Update the selected cell to match the saved view state.
$ws.Range("E30").Select()
